{"js": "// Replace the division-fact answers in the worksheet table.\n// The table has 20 rows x 5 columns; only every 4th row (0, 4, 8, 12, 16)\n// actually holds the division problems, the rows in between are left\n// blank for students to write their own work. We update each of the\n// 25 populated cells in place by (row, col) position so the text is\n// swapped 1:1 without relying on (ambiguous) text search.\nconst replacements = [\n  { row: 0, col: 0, oldText: \"86\u00f79=9, 5\", newText: \"75\u00f73=25, 0\" },\n  { row: 0, col: 1, oldText: \"60\u00f72=30, 0\", newText: \"13\u00f73=4, 1\" },\n  { row: 0, col: 2, oldText: \"17\u00f76=2, 5\", newText: \"63\u00f78=7, 7\" },\n  { row: 0, col: 3, oldText: \"75\u00f76=12, 3\", newText: \"33\u00f72=16, 1\" },\n  { row: 0, col: 4, oldText: \"76\u00f78=9, 4\", newText: \"39\u00f74=9, 3\" },\n\n  { row: 4, col: 0, oldText: \"45\u00f79=5, 0\", newText: \"73\u00f72=36, 1\" },\n  { row: 4, col: 1, oldText: \"92\u00f77=13, 1\", newText: \"76\u00f75=15, 1\" },\n  { row: 4, col: 2, oldText: \"44\u00f74=11, 0\", newText: \"83\u00f73=27, 2\" },\n  { row: 4, col: 3, oldText: \"19\u00f74=4, 3\", newText: \"89\u00f78=11, 1\" },\n  { row: 4, col: 4, oldText: \"93\u00f77=13, 2\", newText: \"43\u00f74=10, 3\" },\n\n  { row: 8, col: 0, oldText: \"18\u00f75=3, 3\", newText: \"75\u00f79=8, 3\" },\n  { row: 8, col: 1, oldText: \"47\u00f74=11, 3\", newText: \"63\u00f76=10, 3\" },\n  { row: 8, col: 2, oldText: \"37\u00f78=4, 5\", newText: \"22\u00f77=3, 1\" },\n  { row: 8, col: 3, oldText: \"61\u00f79=6, 7\", newText: \"60\u00f76=10, 0\" },\n  { row: 8, col: 4, oldText: \"31\u00f73=10, 1\", newText: \"87\u00f78=10, 7\" },\n\n  { row: 12, col: 0, oldText: \"82\u00f73=27, 1\", newText: \"79\u00f74=19, 3\" },\n  { row: 12, col: 1, oldText: \"66\u00f75=13, 1\", newText: \"15\u00f77=2, 1\" },\n  { row: 12, col: 2, oldText: \"12\u00f72=6, 0\", newText: \"37\u00f72=18, 1\" },\n  { row: 12, col: 3, oldText: \"39\u00f74=9, 3\", newText: \"70\u00f72=35, 0\" },\n  { row: 12, col: 4, oldText: \"95\u00f78=11, 7\", newText: \"94\u00f77=13, 3\" },\n\n  { row: 16, col: 0, oldText: \"60\u00f74=15, 0\", newText: \"88\u00f72=44, 0\" },\n  { row: 16, col: 1, oldText: \"94\u00f78=11, 6\", newText: \"98\u00f79=10, 8\" },\n  { row: 16, col: 2, oldText: \"48\u00f78=6, 0\", newText: \"68\u00f74=17, 0\" },\n  { row: 16, col: 3, oldText: \"32\u00f77=4, 4\", newText: \"99\u00f73=33, 0\" },\n  { row: 16, col: 4, oldText: \"93\u00f72=46, 1\", newText: \"79\u00f72=39, 1\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst cells = replacements.map((r) => table.getCell(r.row, r.col));\ncells.forEach((c) => c.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const { newText } = replacements[i];\n  const cell = cells[i];\n  // Setting `.value` replaces only the text run(s) inside the cell's\n  // paragraph(s) while keeping the existing paragraph/run formatting\n  // (font, size, alignment, etc.) untouched.\n  if (cell.value !== newText) {\n    cell.value = newText;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the division-fact answers in the worksheet table.\n#\n# The table is 20 rows x 5 columns; only every 4th row (1, 5, 9, 13, 17 in\n# Word's 1-based indexing) actually holds a division problem -- the rows\n# in between are intentionally blank so students can write their work.\n# We update the 25 populated cells by their (row, col) position, which\n# sidesteps any ambiguity from duplicate/overlapping text values that\n# would make a document-wide Find & Replace unsafe (e.g. one cell's new\n# answer \"39\u00f74=9, 3\" is identical to another cell's original text).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @(\n    @{ Row = 1;  Col = 1; Old = \"86\u00f79=9, 5\";   New = \"75\u00f73=25, 0\" },\n    @{ Row = 1;  Col = 2; Old = \"60\u00f72=30, 0\";  New = \"13\u00f73=4, 1\" },\n    @{ Row = 1;  Col = 3; Old = \"17\u00f76=2, 5\";   New = \"63\u00f78=7, 7\" },\n    @{ Row = 1;  Col = 4; Old = \"75\u00f76=12, 3\";  New = \"33\u00f72=16, 1\" },\n    @{ Row = 1;  Col = 5; Old = \"76\u00f78=9, 4\";   New = \"39\u00f74=9, 3\" },\n\n    @{ Row = 5;  Col = 1; Old = \"45\u00f79=5, 0\";   New = \"73\u00f72=36, 1\" },\n    @{ Row = 5;  Col = 2; Old = \"92\u00f77=13, 1\";  New = \"76\u00f75=15, 1\" },\n    @{ Row = 5;  Col = 3; Old = \"44\u00f74=11, 0\";  New = \"83\u00f73=27, 2\" },\n    @{ Row = 5;  Col = 4; Old = \"19\u00f74=4, 3\";   New = \"89\u00f78=11, 1\" },\n    @{ Row = 5;  Col = 5; Old = \"93\u00f77=13, 2\";  New = \"43\u00f74=10, 3\" },\n\n    @{ Row = 9;  Col = 1; Old = \"18\u00f75=3, 3\";   New = \"75\u00f79=8, 3\" },\n    @{ Row = 9;  Col = 2; Old = \"47\u00f74=11, 3\";  New = \"63\u00f76=10, 3\" },\n    @{ Row = 9;  Col = 3; Old = \"37\u00f78=4, 5\";   New = \"22\u00f77=3, 1\" },\n    @{ Row = 9;  Col = 4; Old = \"61\u00f79=6, 7\";   New = \"60\u00f76=10, 0\" },\n    @{ Row = 9;  Col = 5; Old = \"31\u00f73=10, 1\";  New = \"87\u00f78=10, 7\" },\n\n    @{ Row = 13; Col = 1; Old = \"82\u00f73=27, 1\";  New = \"79\u00f74=19, 3\" },\n    @{ Row = 13; Col = 2; Old = \"66\u00f75=13, 1\";  New = \"15\u00f77=2, 1\" },\n    @{ Row = 13; Col = 3; Old = \"12\u00f72=6, 0\";   New = \"37\u00f72=18, 1\" },\n    @{ Row = 13; Col = 4; Old = \"39\u00f74=9, 3\";   New = \"70\u00f72=35, 0\" },\n    @{ Row = 13; Col = 5; Old = \"95\u00f78=11, 7\";  New = \"94\u00f77=13, 3\" },\n\n    @{ Row = 17; Col = 1; Old = \"60\u00f74=15, 0\";  New = \"88\u00f72=44, 0\" },\n    @{ Row = 17; Col = 2; Old = \"94\u00f78=11, 6\";  New = \"98\u00f79=10, 8\" },\n    @{ Row = 17; Col = 3; Old = \"48\u00f78=6, 0\";   New = \"68\u00f74=17, 0\" },\n    @{ Row = 17; Col = 4; Old = \"32\u00f77=4, 4\";   New = \"99\u00f73=33, 0\" },\n    @{ Row = 17; Col = 5; Old = \"93\u00f72=46, 1\";  New = \"79\u00f72=39, 1\" }\n)\n\nforeach ($u in $updates) {\n    $cell = $t.Cell($u.Row, $u.Col)\n    # Assigning Range.Text replaces the cell's text while leaving the\n    # paragraph/run formatting (font, size, alignment) untouched.\n    $cell.Range.Text = $u.New\n}\n"}
